$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.880.85'
$ws.Range("E2").Value = '  -2.68%  '
$ws.Range("D3").Value = '1.808.88'
$ws.Range("E3").Value = '  -3.55%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '231.53'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("E6").Value = '  -1.65%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("D8").Value = '39.05'
$ws.Range("E8").Value = '  -7.65%  '
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").Value = '0.0678'
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("E11").Value = '  -1.87%  '
$ws.Range("D12").Value = '2.070.01'
$ws.Range("E12").Value = '  -3.59%  '
$ws.Range("D13").Value = '1.812.05'
$ws.Range("E13").Value = '  -3.47%  '
$ws.Range("D14").Value = '0.659'
$ws.Range("E14").Value = '  -4.05%  '
$ws.Range("E15").Value = '  -6.90%  '
$ws.Range("D16").Value = '4.55'
$ws.Range("D17").Value = '34.814.80'
$ws.Range("E17").Value = '  -2.99%  '
$ws.Range("D18").Value = '68.97'
$ws.Range("E18").Value = '  -2.60%  '
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("D20").Value = '237.83'
$ws.Range("E20").Value = '  -4.13%  '
$ws.Range("D21").Value = '11.76'
$ws.Range("E21").Value = '  -6.19%  '
$ws.Range("D22").Value = '4.62'
$ws.Range("E22").Value = '  -4.16%  '
$ws.Range("E23").Value = '  +0.27%  '
$ws.Range("E24").Value = '  -0.54%  '
$ws.Range("D25").Value = '171.84'
$ws.Range("E25").Value = '  +0.23%  '
$ws.Range("E26").Value = '  -4.29%  '
$ws.Range("D27").Value = '17.18'
$ws.Range("E27").Value = '  -4.79%  '
$ws.Range("E28").Value = '  -3.55%  '
$ws.Range("D29").Value = '1.54'
$ws.Range("E29").Value = '  +5.94%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").Value = '  -0.89%  '
$ws.Range("E34").Value = '  -8.88%  '
$ws.Range("D35").Value = '1.15'
$ws.Range("E35").Value = '  +4.97%  '
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("D37").Value = '91.06'
$ws.Range("E37").Value = '  -7.76%  '
$ws.Range("D38").Value = '1.31'
$ws.Range("E38").Value = '  +3.17%  '
$ws.Range("D39").Value = '1.304.60'
$ws.Range("E39").Value = '  -4.56%  '
$ws.Range("E40").Value = '  -3.60%  '
$ws.Range("D41").Value = '2.45'
$ws.Range("E41").Value = '  -1.51%  '
$ws.Range("D42").Value = '0.958'
$ws.Range("E42").Value = '  -7.22%  '
$ws.Range("D43").Value = '14.31'
$ws.Range("E43").Value = '  -5.33%  '
$ws.Range("E44").Value = '  -13.20%  '
$ws.Range("E45").Value = '  -5.27%  '
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("E47").Value = '  -2.15%  '
$ws.Range("D48").Value = '1.991.41'
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("E49").Value = '  +7.24%  '
$ws.Range("E50").Value = '  +0.29%  '
$ws.Range("D51").Value = '98.69'
$ws.Range("E51").Value = '  -6.48%  '
